$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.66"
$ws.Range("E2").Value = "'0.93%"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'27.26"
$ws.Range("E3").Value = "'2.20%"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'4.860"
$ws.Range("E4").Value = "'-0.38%"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.06403"
$ws.Range("E5").Value = "'1.47%"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'7.016"
$ws.Range("G6").Value = "'12"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'1.194"
$ws.Range("E7").Value = "'-8.26%"
$ws.Range("G7").Value = "'12"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.8868"
$ws.Range("E8").Value = "'1.95%"
$ws.Range("G8").Value = "'12"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1540"
$ws.Range("E9").Value = "'0.85%"
$ws.Range("G9").Value = "'12"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.05170"
$ws.Range("E10").Value = "'1.93%"
$ws.Range("G10").Value = "'12"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07505"
$ws.Range("E11").Value = "'1.10%"
$ws.Range("G11").Value = "'12"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02884"
$ws.Range("E12").Value = "'-0.59%"
$ws.Range("G12").Value = "'12"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.08968"
$ws.Range("E13").Value = "'-1.00%"
$ws.Range("G13").Value = "'12"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001573"
$ws.Range("E14").Value = "'-0.36%"
$ws.Range("G14").Value = "'12"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0006351"
$ws.Range("E15").Value = "'0.30%"
$ws.Range("G15").Value = "'12"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.006127"
$ws.Range("E16").Value = "'2.22%"
$ws.Range("G16").Value = "'12"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.475"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("G17").Value = "'12"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = "'3.304"
$ws.Range("E18").Value = "'-0.24%"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'2.247"
$ws.Range("E19").Value = "'-1.61%"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.3085"
$ws.Range("E20").Value = "'-2.03%"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.1342"
$ws.Range("E21").Value = "'2.21%"
$ws.Range("G21").Value = "'12"
$ws.Range("D22").Value = "'3.910"
$ws.Range("E22").Value = "'0.09%"
$ws.Range("G22").Value = "'12"
$ws.Range("E23").Value = "'10.02%"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'0.04423"
$ws.Range("E24").Value = "'0.84%"
$ws.Range("G24").Value = "'12"
$ws.Range("E25").Value = "'0.45%"
$ws.Range("G25").Value = "'12"
$ws.Range("D26").Value = "'0.003894"
$ws.Range("E26").Value = "'-7.50%"
$ws.Range("G26").Value = "'12"
$ws.Range("G27").Value = "'12"
$ws.Range("E28").Value = "'-1.82%"
$ws.Range("G28").Value = "'12"
$ws.Range("D29").Value = "'0.0001643"
$ws.Range("E29").Value = "'-1.04%"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.04108"
$ws.Range("E40").Value = "'0.87%"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.006809"
$ws.Range("E41").Value = "'-3.42%"
$ws.Range("G41").Value = "'12"
$ws.Range("E42").Value = "'0.34%"
$ws.Range("G42").Value = "'12"
$ws.Range("E43").Value = "'-12.68%"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.01166"
$ws.Range("E44").Value = "'4.50%"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005336"
$ws.Range("E45").Value = "'2.54%"
$ws.Range("G45").Value = "'12"
$ws.Range("D46").Value = "'1.561"
$ws.Range("E46").Value = "'5.01%"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.01850"
$ws.Range("E47").Value = "'-7.47%"
$ws.Range("G47").Value = "'12"
$ws.Range("G48").Value = "'12"
$ws.Range("G49").Value = "'12"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
